$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("individual")

# Row 7 (mediacategory_media): update the display name to drop "mass"
$ws.Range("B7").Value = "Gets COVID-19 related information and advice from traditional media outlets"

# Row 6 (income): correct the filter/ordering to be proper JSON object mappings
$ws.Range("D6").Value = '{"i prefer not to say":13}'
$ws.Range("E6").Value = '{"<€500":1, "€501-€1000":2, "€1001-€1500":3, "€1501-€2000":4, "€2001-€2500":5, "€2501-€3000":6, "€3001-€3500":7, "€3501-€4000":8, "€4001-€4500":9, "€4501-€5000":10, "€5001-€7500":11, ">€7500":12}'

# Update selection to reflect the edited cell
$ws.Range("E6").Select()
